# Update test case sheets for registration, search, vaccine administer, statistic
# (checkmarks added to the "Passed or not" column; Vaccine Administration sheet's
#  empty "Actual Result" column removed; active tab moved to Patient Search Module)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 2 - "Vaccine Administration Module"
# Drop the unused "Actual Result" column (column B) so the old column C
# ("Passed or not") becomes the new column B, then mark every scenario row
# as passed with a check mark.
# ---------------------------------------------------------------------------
$wsVaccine = $wb.Worksheets.Item("Vaccine Administration Module")
$wsVaccine.Columns.Item(2).Delete() | Out-Null

$wsVaccine.Range("B3").Value  = "√"
$wsVaccine.Range("B5").Value  = "√"
$wsVaccine.Range("B7").Value  = "√"
$wsVaccine.Range("B11").Value = "√"
$wsVaccine.Range("B13").Value = "√"
$wsVaccine.Range("B15").Value = "√"
$wsVaccine.Range("B17").Value = "√"
$wsVaccine.Range("B19").Value = "√"

$wsVaccine.Range("B11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3 - "Patient Search Module"
# Mark every scenario row as passed with a check mark.
# ---------------------------------------------------------------------------
$wsSearch = $wb.Worksheets.Item("Patient Search Module")

$wsSearch.Range("B3").Value  = "√"
$wsSearch.Range("B5").Value  = "√"
$wsSearch.Range("B7").Value  = "√"
$wsSearch.Range("B11").Value = "√"
$wsSearch.Range("B13").Value = "√"

$wsSearch.Range("B23").Select() | Out-Null

# Patient Search Module becomes the active / displayed sheet.
$wsSearch.Activate() | Out-Null
